# faturamento_diario.xlsx update
# Insert a new daily-revenue row for day 14 of August/2025 right before the
# existing row 15 (first row of the July/2025 block), shifting the July/June/
# May rows down by one. This matches the source workbook which now also
# has an extra day of data for August.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 15:106 down to 16:107 to make room for the new row.
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the August 14th figures.
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = 21233.96
$ws.Cells.Item(15, 3).Value = 8
$ws.Cells.Item(15, 4).Value = 2025
$ws.Cells.Item(15, 5).Value = "08/2025"
